$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.278710430081162
$ws.Cells.Item(2, 3).Value = 0.3360007219306169
$ws.Cells.Item(2, 4).Value = 0.6466354191979633
$ws.Cells.Item(2, 5).Value = 0.2641277278156693
$ws.Cells.Item(2, 7).Value = 0.5101622322586152
$ws.Cells.Item(2, 8).Value = 0.6506829541225727
$ws.Cells.Item(2, 9).Value = 0.4235443996883959
$ws.Cells.Item(2, 10).Value = 0.1377940913252402
$ws.Cells.Item(2, 15).Value = 2.287274129948258

$ws.Cells.Item(3, 2).Value = 1.133085885519336
$ws.Cells.Item(3, 3).Value = 0.2944377662282136
$ws.Cells.Item(3, 4).Value = 0.6358277581243783
$ws.Cells.Item(3, 5).Value = 0.2588352462145806
$ws.Cells.Item(3, 7).Value = 0.5137103184422784
$ws.Cells.Item(3, 8).Value = 0.6575086674885071
$ws.Cells.Item(3, 9).Value = 0.4336482092987435
$ws.Cells.Item(3, 10).Value = 0.1342939655858189
$ws.Cells.Item(3, 15).Value = 2.308749860258999

$ws.Cells.Item(4, 2).Value = 1.043426673278077
$ws.Cells.Item(4, 3).Value = 0.2688236922782608
$ws.Cells.Item(4, 4).Value = 0.6295269056035124
$ws.Cells.Item(4, 5).Value = 0.2557286836162547
$ws.Cells.Item(4, 7).Value = 0.516423489300827
$ws.Cells.Item(4, 8).Value = 0.6621207783661802
$ws.Cells.Item(4, 9).Value = 0.4402598989969366
$ws.Cells.Item(4, 10).Value = 0.1322237867784892
$ws.Cells.Item(4, 15).Value = 2.323939402934784

$ws.Cells.Item(5, 2).Value = 1.00683058098457
$ws.Cells.Item(5, 3).Value = 0.258362724150544
$ws.Cells.Item(5, 4).Value = 0.6270435832576027
$ws.Cells.Item(5, 5).Value = 0.2544987092132587
$ws.Cells.Item(5, 7).Value = 0.5176630958995716
$ws.Cells.Item(5, 8).Value = 0.6641060250547
$ws.Cells.Item(5, 9).Value = 0.4430565132271802
$ws.Cells.Item(5, 10).Value = 0.1313999999468365
$ws.Cells.Item(5, 15).Value = 2.330631783171143

$ws.Cells.Item(6, 2).Value = 1.000750321524208
$ws.Cells.Item(6, 3).Value = 0.2566243177181207
$ws.Cells.Item(6, 4).Value = 0.626636324492182
$ws.Cells.Item(6, 5).Value = 0.2542966459020874
$ws.Cells.Item(6, 7).Value = 0.5178770091465665
$ws.Cells.Item(6, 8).Value = 0.6644420589170252
$ws.Cells.Item(6, 9).Value = 0.4435270593059109
$ws.Cells.Item(6, 10).Value = 0.1312644078867464
$ws.Cells.Item(6, 15).Value = 2.331773362249393

$ws.Cells.Item(7, 2).Value = 1.04293336220735
$ws.Cells.Item(7, 3).Value = 0.2686827043293931
$ws.Cells.Item(7, 4).Value = 0.629493073130277
$ws.Cells.Item(7, 5).Value = 0.2557119500808369
$ws.Cells.Item(7, 7).Value = 0.5164396652651888
$ws.Cells.Item(7, 8).Value = 0.6621471239750107
$ws.Cells.Item(7, 9).Value = 0.4402972012842188
$ws.Cells.Item(7, 10).Value = 0.1322125966386309
$ws.Cells.Item(7, 15).Value = 2.324027625693702

$ws.Cells.Item(8, 2).Value = 1.228551362871031
$ws.Cells.Item(8, 3).Value = 0.3216897671637753
$ws.Cells.Item(8, 4).Value = 0.6428394216509616
$ws.Cells.Item(8, 5).Value = 0.2622731927435993
$ws.Cells.Item(8, 7).Value = 0.5112743662953037
$ws.Cells.Item(8, 8).Value = 0.6529489933594306
$ws.Cells.Item(8, 9).Value = 0.4269433673916421
$ws.Cells.Item(8, 10).Value = 0.1365708589184962
$ws.Cells.Item(8, 15).Value = 2.294262416528866

$ws.Cells.Item(9, 2).Value = 1.590514941090589
$ws.Cells.Item(9, 3).Value = 0.4248641004353431
$ws.Cells.Item(9, 4).Value = 0.6716698084565564
$ws.Cells.Item(9, 5).Value = 0.2762755261495258
$ws.Cells.Item(9, 7).Value = 0.5054089239649358
$ws.Cells.Item(9, 8).Value = 0.638257921127078
$ws.Cells.Item(9, 9).Value = 0.4040050565492361
$ws.Cells.Item(9, 10).Value = 0.1457448596886266
$ws.Cells.Item(9, 15).Value = 2.251845834517695

$ws.Cells.Item(10, 2).Value = 1.855119626744397
$ws.Cells.Item(10, 3).Value = 0.5001700646299128
$ws.Cells.Item(10, 4).Value = 0.6944745126083944
$ws.Cells.Item(10, 5).Value = 0.2872580120580395
$ws.Cells.Item(10, 7).Value = 0.5037294056500485
$ws.Cells.Item(10, 8).Value = 0.6295118566182083
$ws.Cells.Item(10, 9).Value = 0.3891484936105272
$ws.Cells.Item(10, 10).Value = 0.1528702646424733
$ws.Cells.Item(10, 15).Value = 2.230488437656135

$ws.Cells.Item(11, 2).Value = 1.975189165056747
$ws.Cells.Item(11, 3).Value = 0.534316136389009
$ws.Cells.Item(11, 4).Value = 0.7052021442258933
$ws.Cells.Item(11, 5).Value = 0.2924057916915572
$ws.Cells.Item(11, 7).Value = 0.5035429848869626
$ws.Cells.Item(11, 8).Value = 0.6259792461880664
$ws.Cells.Item(11, 9).Value = 0.3828266901222239
$ws.Cells.Item(11, 10).Value = 0.1561961037373862
$ws.Cells.Item(11, 15).Value = 2.22291916028135

$ws.Cells.Item(12, 2).Value = 2.020611182450295
$ws.Cells.Item(12, 3).Value = 0.5472298471112254
$ws.Cells.Item(12, 4).Value = 0.7093152696487834
$ws.Cells.Item(12, 5).Value = 0.2943769716818068
$ws.Cells.Item(12, 7).Value = 0.5035559716274207
$ws.Cells.Item(12, 8).Value = 0.624705809559984
$ws.Cells.Item(12, 9).Value = 0.3804958528995481
$ws.Cells.Item(12, 10).Value = 0.1574676908864774
$ws.Cells.Item(12, 15).Value = 2.220362902176248

$ws.Cells.Item(13, 2).Value = 2.010830810894788
$ws.Cells.Item(13, 3).Value = 0.5444494008106062
$ws.Cells.Item(13, 4).Value = 0.7084271757760519
$ws.Cells.Item(13, 5).Value = 0.2939514720096241
$ws.Cells.Item(13, 7).Value = 0.5035494491494887
$ws.Cells.Item(13, 8).Value = 0.6249772053924261
$ws.Cells.Item(13, 9).Value = 0.3809950292250068
$ws.Cells.Item(13, 10).Value = 0.1571932905327742
$ws.Cells.Item(13, 15).Value = 2.220899625279429

$ws.Cells.Item(14, 2).Value = 1.978926992129118
$ws.Cells.Item(14, 3).Value = 0.5353788931584518
$ws.Cells.Item(14, 4).Value = 0.7055395156302495
$ws.Cells.Item(14, 5).Value = 0.2925675243967873
$ws.Cells.Item(14, 7).Value = 0.5035423752565578
$ws.Cells.Item(14, 8).Value = 0.6258731905414265
$ws.Cells.Item(14, 9).Value = 0.3826336644320918
$ws.Cells.Item(14, 10).Value = 0.1563004740802967
$ws.Cells.Item(14, 15).Value = 2.22270263310557

$ws.Cells.Item(15, 2).Value = 1.959378955147827
$ws.Cells.Item(15, 3).Value = 0.529820754156674
$ws.Cells.Item(15, 4).Value = 0.7037773551032274
$ws.Cells.Item(15, 5).Value = 0.2917226591261226
$ws.Cells.Item(15, 7).Value = 0.5035489423750619
$ws.Cells.Item(15, 8).Value = 0.6264303837224787
$ws.Cells.Item(15, 9).Value = 0.3836456018593175
$ws.Cells.Item(15, 10).Value = 0.1557551829204584
$ws.Cells.Item(15, 15).Value = 2.22384744909121

$ws.Cells.Item(16, 2).Value = 1.847266539228769
$ws.Cells.Item(16, 3).Value = 0.4979362446733262
$ws.Cells.Item(16, 4).Value = 0.6937805487325193
$ws.Cells.Item(16, 5).Value = 0.2869246469899736
$ws.Cells.Item(16, 7).Value = 0.5037532561372728
$ws.Cells.Item(16, 8).Value = 0.6297517088129041
$ws.Cells.Item(16, 9).Value = 0.3895704497309254
$ws.Cells.Item(16, 10).Value = 0.1526546143647494
$ws.Cells.Item(16, 15).Value = 2.231026418954087

$ws.Cells.Item(17, 2).Value = 1.778410456996767
$ws.Cells.Item(17, 3).Value = 0.4783471995091872
$ws.Cells.Item(17, 4).Value = 0.687738375127509
$ws.Cells.Item(17, 5).Value = 0.2840200994031008
$ws.Cells.Item(17, 7).Value = 0.5040269052134363
$ws.Cells.Item(17, 8).Value = 0.6319035764133929
$ws.Cells.Item(17, 9).Value = 0.3933171842487706
$ws.Cells.Item(17, 10).Value = 0.1507741610209337
$ws.Cells.Item(17, 15).Value = 2.235981221530011

$ws.Cells.Item(18, 2).Value = 1.738778126882039
$ws.Cells.Item(18, 3).Value = 0.4670696937232606
$ws.Cells.Item(18, 4).Value = 0.6842963675974829
$ws.Cells.Item(18, 5).Value = 0.2823637677812343
$ws.Cells.Item(18, 7).Value = 0.5042386363012099
$ws.Cells.Item(18, 8).Value = 0.6331832461406179
$ws.Cells.Item(18, 9).Value = 0.3955132751063939
$ws.Cells.Item(18, 10).Value = 0.1497005222884127
$ws.Cells.Item(18, 15).Value = 2.239033028631582

$ws.Cells.Item(19, 2).Value = 1.725354545838456
$ws.Cells.Item(19, 3).Value = 0.4632495585693164
$ws.Cells.Item(19, 4).Value = 0.6831366823147675
$ws.Cells.Item(19, 5).Value = 0.2818054162185035
$ws.Cells.Item(19, 7).Value = 0.5043196410378386
$ws.Cells.Item(19, 8).Value = 0.6336237245749174
$ws.Cells.Item(19, 9).Value = 0.3962638774620313
$ws.Cells.Item(19, 10).Value = 0.1493383712033136
$ws.Cells.Item(19, 15).Value = 2.240100957239889

$ws.Cells.Item(20, 2).Value = 1.785743234370727
$ws.Cells.Item(20, 3).Value = 0.4804335688355081
$ws.Cells.Item(20, 4).Value = 0.6883781297132998
$ws.Cells.Item(20, 5).Value = 0.2843278147686803
$ws.Cells.Item(20, 7).Value = 0.5039921475076596
$ws.Cells.Item(20, 8).Value = 0.6316701613639992
$ws.Cells.Item(20, 9).Value = 0.3929140852230439
$ws.Cells.Item(20, 10).Value = 0.1509735157323036
$ws.Cells.Item(20, 15).Value = 2.235432864400678

$ws.Cells.Item(21, 2).Value = 1.988299177467809
$ws.Cells.Item(21, 3).Value = 0.5380435790145839
$ws.Cells.Item(21, 4).Value = 0.7063863130082382
$ws.Cells.Item(21, 5).Value = 0.2929734309616592
$ws.Cells.Item(21, 7).Value = 0.5035421805336853
$ws.Cells.Item(21, 8).Value = 0.625608272082232
$ws.Cells.Item(21, 9).Value = 0.3821506426763115
$ws.Cells.Item(21, 10).Value = 0.156562385630437
$ws.Cells.Item(21, 15).Value = 2.222164619628217

$ws.Cells.Item(22, 2).Value = 2.120413614560618
$ws.Cells.Item(22, 3).Value = 0.5755975891067919
$ws.Cells.Item(22, 4).Value = 0.7184518040891703
$ws.Cells.Item(22, 5).Value = 0.2987510689998203
$ws.Cells.Item(22, 7).Value = 0.5037355119726357
$ws.Cells.Item(22, 8).Value = 0.6220212560223786
$ws.Cells.Item(22, 9).Value = 0.3754840115506628
$ws.Cells.Item(22, 10).Value = 0.1602859457748878
$ws.Cells.Item(22, 15).Value = 2.215300998866326

$ws.Cells.Item(23, 2).Value = 2.049926957289415
$ws.Cells.Item(23, 3).Value = 0.5555634560772091
$ws.Cells.Item(23, 4).Value = 0.7119851477734755
$ws.Cells.Item(23, 5).Value = 0.295655793332557
$ws.Cells.Item(23, 7).Value = 0.5035875571034865
$ws.Cells.Item(23, 8).Value = 0.6239013723092199
$ws.Cells.Item(23, 9).Value = 0.3790083489523095
$ws.Cells.Item(23, 10).Value = 0.1582921168900242
$ws.Cells.Item(23, 15).Value = 2.218798344890246

$ws.Cells.Item(24, 2).Value = 1.782428228976073
$ws.Cells.Item(24, 3).Value = 0.4794903695030257
$ws.Cells.Item(24, 4).Value = 0.6880887978635144
$ws.Cells.Item(24, 5).Value = 0.2841886545043977
$ws.Cells.Item(24, 7).Value = 0.5040076920311662
$ws.Cells.Item(24, 8).Value = 0.6317755558302451
$ws.Cells.Item(24, 9).Value = 0.3930961953522356
$ws.Cells.Item(24, 10).Value = 0.1508833642363498
$ws.Cells.Item(24, 15).Value = 2.23568014372259

$ws.Cells.Item(25, 2).Value = 1.492821311943601
$ws.Cells.Item(25, 3).Value = 0.3970380475581123
$ws.Cells.Item(25, 4).Value = 0.663585588990486
$ws.Cells.Item(25, 5).Value = 0.2723656678298525
$ws.Cells.Item(25, 7).Value = 0.5065360600915199
$ws.Cells.Item(25, 8).Value = 0.6418732232329489
$ws.Cells.Item(25, 9).Value = 0.4098610651457264
$ws.Cells.Item(25, 10).Value = 0.1431955919822414
$ws.Cells.Item(25, 15).Value = 2.261604419784135
